# Generate Report for Handoff
# Updates the localization-status workbook so that the "b.md" rows reflect
# that a new handoff package has been generated for it (zh-cn and de-de),
# instead of showing the same "handed back" status as a.md.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is for b.md. Its zh-cn / de-de status columns and
# "Latest HO Xliff Generate Date" column need to move from the old
# "Handed back: in sync with en-US" / 2016-09-04 18:41:20 values to the
# new "Ready for handoff" / 2016-09-04 18:42:10 values.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 18:42:10"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) gets a fresh handoff file/status.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-04 18:42:03"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16289712e6e111d66615baddb73fdbc2b06a989f/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c0a3677022657ff5ffa5b60b7ef689a9d8d094ae/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: row 3 (b.md) gets a fresh handoff file/status.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-04 18:42:10"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16289712e6e111d66615baddb73fdbc2b06a989f/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c0a3677022657ff5ffa5b60b7ef689a9d8d094ae/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17

$wb.Save()
